$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old schedule entirely (dates, topics, styles, etc.)
$ws.Range("A1:F20").Clear()

# --- Dates (column A), formatted like the original "d-mmm" style ---
$dates = @(43683,43684,43685,43686,43687,43688,43689,43690,43691,43692,43693,43694,43695,43696)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 1
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "d-mmm"
}

# --- Topic / note text (columns C and D) ---
$ws.Range("C1").Value = "Carbonyls and Alcohols"
$ws.Range("D1").Value = "CARS"

$ws.Range("C2").Value = "Expression of Genetic Information"
$ws.Range("D2").Value = "Carbohydrates"

$ws.Range("C3").Value = "Metabolic Components"
$ws.Range("D3").Value = "Metabolic Pathways"

$ws.Range("C4").Value = "Nitrogen Chemistry"
$ws.Range("D4").Value = "Catchup"

$ws.Range("C5").Value = "Review all notes. Print all notes."

$ws.Range("C6").Value = "Practice Exam 1"

$ws.Range("C13").Value = "Practice Exam 2"

$ws.Range("C7").Value = "Review practice exam. Do practice problems. Study flashcards."

# --- Wrap + left-align the merged "review" block (C7:C12), then propagate the
#     formatting without minting duplicate style records ---
$c7 = $ws.Range("C7")
$c7.WrapText = $true
$c7.HorizontalAlignment = -4131
$c7.Copy()
$ws.Range("C8:C12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Merge the review block ---
$ws.Range("C7:C12").Merge()

# --- Column C width ---
$ws.Columns.Item(3).ColumnWidth = 28.666666666667

# --- Selection marker ---
$ws.Range("E10").Select()
